# Update the acquisition-timestamp column (A) for the existing rows on the
# "ランサーズ" sheet: the scraper re-ran at 2025-10-10 18:25:59 and the
# previously-recorded timestamp (2025-10-10 12:46:22) for rows 2-17 is
# refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-10 12:46:22"
$newTimestamp = "2025-10-10 18:25:59"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
